# Fix checkout: DB-only (no Excel) + keep mp.initPoint response
# Set the "stock" column (C) to 10 for every card row, and update the
# active selection to C9 (matching the author's last selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stock values (column C, rows 2-12) to 10
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 10
}

# Move the selection to C9 as recorded in the saved view state
$ws.Range("C9").Select()
